# Auto-generated Excel COM-interop edit script
# Applies the numeric corrections from the commit diff to the
# eight leve-profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H57").Value = 99989.5
$ws.Range("J57").Value = 99989.5
$ws.Range("L57").Value = 299968.5
$ws.Range("N57").Value = -300966.5

$ws.Range("H116").Value = 364095.1
$ws.Range("I116").Value = 11750
$ws.Range("J116").Value = 452181.38
$ws.Range("K116").Value = 11750
$ws.Range("L116").Value = 452181.38
$ws.Range("M116").Value = -8308
$ws.Range("N116").Value = -459065.38

$ws.Range("H132").Value = 76029.09
$ws.Range("I132").Value = 88816.11
$ws.Range("K132").Value = 266448.33
$ws.Range("M132").Value = -263918.33

$ws.Range("H133").Value = 99567.71000000001
$ws.Range("J133").Value = 99567.71000000001
$ws.Range("L133").Value = 99567.71000000001
$ws.Range("N133").Value = -109687.71

$ws.Range("H136").Value = 83999.25
$ws.Range("J136").Value = 83999.25
$ws.Range("L136").Value = 83999.25
$ws.Range("N136").Value = -94199.25

$ws.Range("H139").Value = 94998
$ws.Range("J139").Value = 94998
$ws.Range("L139").Value = 94998
$ws.Range("N139").Value = -105278

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5265429.5
$ws.Range("I32").Value = 6251421
$ws.Range("K32").Value = 6251421
$ws.Range("M32").Value = -6251134

$ws.Range("H74").Value = 3790319.8
$ws.Range("I74").Value = 4809255
$ws.Range("K74").Value = 4809255
$ws.Range("M74").Value = -4808381

$ws.Range("H76").Value = 28420
$ws.Range("J76").Value = 37499.5
$ws.Range("L76").Value = 37499.5
$ws.Range("N76").Value = -38175.5

$ws.Range("H77").Value = 3790319.8
$ws.Range("I77").Value = 4809255
$ws.Range("K77").Value = 24046275
$ws.Range("M77").Value = -24041907

$ws.Range("H79").Value = 28420
$ws.Range("J79").Value = 37499.5
$ws.Range("L79").Value = 37499.5
$ws.Range("N79").Value = -39839.5

$ws.Range("H102").Value = 43199.668
$ws.Range("I102").Value = 48099.625
$ws.Range("K102").Value = 48099.625
$ws.Range("M102").Value = -46477.625

$ws.Range("H121").Value = 79999.5
$ws.Range("J121").Value = 79999.5
$ws.Range("L121").Value = 79999.5
$ws.Range("N121").Value = -83493.5

$ws.Range("H123").Value = 74999
$ws.Range("J123").Value = 74999
$ws.Range("L123").Value = 74999
$ws.Range("N123").Value = -84799

$ws.Range("H125").Value = 75999
$ws.Range("J125").Value = 75999
$ws.Range("L125").Value = 75999
$ws.Range("N125").Value = -85839

$ws.Range("H132").Value = 18772088
$ws.Range("I132").Value = 23776776
$ws.Range("J132").Value = 4512.25
$ws.Range("K132").Value = 71330328
$ws.Range("L132").Value = 13536.75
$ws.Range("M132").Value = -71327798
$ws.Range("N132").Value = -18596.75

$ws.Range("H134").Value = 166713710
$ws.Range("J134").Value = 166713710
$ws.Range("L134").Value = 166713710
$ws.Range("N134").Value = -166723850

$ws.Range("H137").Value = 1000000000
$ws.Range("J137").Value = 1000000000
$ws.Range("L137").Value = 1000000000
$ws.Range("N137").Value = -1000010200

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H140").Value = 83965.60000000001
$ws.Range("J140").Value = 89699.336
$ws.Range("L140").Value = 89699.336
$ws.Range("N140").Value = -100059.336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 999.7059
$ws.Range("I22").Value = 998.3333
$ws.Range("K22").Value = 998.3333
$ws.Range("M22").Value = -648.3333

$ws.Range("H31").Value = 72751.375
$ws.Range("J31").Value = 37654.547
$ws.Range("L31").Value = 37654.547
$ws.Range("N31").Value = -38244.547

$ws.Range("H34").Value = 72751.375
$ws.Range("J34").Value = 37654.547
$ws.Range("L34").Value = 37654.547
$ws.Range("N34").Value = -38058.547

$ws.Range("H58").Value = 6524216
$ws.Range("I58").Value = 8864819
$ws.Range("K58").Value = 8864819
$ws.Range("M58").Value = -8864616

$ws.Range("H68").Value = 97236
$ws.Range("J68").Value = 97236
$ws.Range("L68").Value = 97236
$ws.Range("N68").Value = -98734

$ws.Range("H71").Value = 97236
$ws.Range("J71").Value = 97236
$ws.Range("L71").Value = 291708
$ws.Range("N71").Value = -299196

$ws.Range("H132").Value = 79185990
$ws.Range("I132").Value = 83356980
$ws.Range("J132").Value = 62502056
$ws.Range("K132").Value = 250070940
$ws.Range("L132").Value = 187506168
$ws.Range("M132").Value = -250068410
$ws.Range("N132").Value = -187511228

$ws.Range("H134").Value = 41708612
$ws.Range("I134").Value = 41708612
$ws.Range("K134").Value = 125125836
$ws.Range("M134").Value = -125123301

$ws.Range("H136").Value = 6524216
$ws.Range("I136").Value = 8864819
$ws.Range("K136").Value = 26594457
$ws.Range("M136").Value = -26591907

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 99325.75
$ws.Range("J37").Value = 99325.75
$ws.Range("L37").Value = 297977.25
$ws.Range("N37").Value = -298201.25

$ws.Range("H98").Value = 497.33334
$ws.Range("J98").Value = 292
$ws.Range("L98").Value = 876
$ws.Range("N98").Value = -3872

$ws.Range("H132").Value = 1599.6
$ws.Range("I132").Value = 1334.3334
$ws.Range("K132").Value = 12009.0006
$ws.Range("M132").Value = -9479.000599999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 7114.143
$ws.Range("J43").Value = 12000
$ws.Range("L43").Value = 12000
$ws.Range("N43").Value = -12302

$ws.Range("H122").Value = 5681.4
$ws.Range("I122").Value = 3429.1765
$ws.Range("J122").Value = 10467.375
$ws.Range("K122").Value = 10287.5295
$ws.Range("L122").Value = 31402.125
$ws.Range("M122").Value = -7837.529500000001
$ws.Range("N122").Value = -36302.125

$ws.Range("H123").Value = 16333
$ws.Range("J123").Value = 16333
$ws.Range("L123").Value = 16333
$ws.Range("N123").Value = -21233

$ws.Range("H126").Value = 1393938.2
$ws.Range("J126").Value = 4800
$ws.Range("L126").Value = 14400
$ws.Range("N126").Value = -19340

$ws.Range("H132").Value = 384152.84
$ws.Range("I132").Value = 505182.84
$ws.Range("K132").Value = 1515548.52
$ws.Range("M132").Value = -1513018.52

$ws.Range("H134").Value = 49000
$ws.Range("J134").Value = 49000
$ws.Range("L134").Value = 147000
$ws.Range("N134").Value = -152070

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 959.1177
$ws.Range("I22").Value = 959.0345
$ws.Range("J22").Value = 959.6
$ws.Range("K22").Value = 959.0345
$ws.Range("L22").Value = 959.6
$ws.Range("M22").Value = -664.0345
$ws.Range("N22").Value = -1549.6

$ws.Range("H27").Value = 959.1177
$ws.Range("I27").Value = 959.0345
$ws.Range("J27").Value = 959.6
$ws.Range("K27").Value = 959.0345
$ws.Range("L27").Value = 959.6
$ws.Range("M27").Value = -852.0345
$ws.Range("N27").Value = -1173.6

$ws.Range("H46").Value = 14676.4
$ws.Range("J46").Value = 1000
$ws.Range("L46").Value = 1000
$ws.Range("N46").Value = -1376

$ws.Range("H132").Value = 16736204
$ws.Range("I132").Value = 17532976
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 52598928
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -52596398
$ws.Range("N132").Value = -17060

$ws.Range("H136").Value = 12871018
$ws.Range("J136").Value = 420410.34
$ws.Range("L136").Value = 1261231.02
$ws.Range("N136").Value = -1266331.02

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2424.0557
$ws.Range("I122").Value = 1792.1428
$ws.Range("K122").Value = 5376.428400000001
$ws.Range("M122").Value = -2926.428400000001

$ws.Range("H132").Value = 14993577
$ws.Range("I132").Value = 8918354
$ws.Range("J132").Value = 83339830
$ws.Range("K132").Value = 26755062
$ws.Range("L132").Value = 250019490
$ws.Range("M132").Value = -26752532
$ws.Range("N132").Value = -250024550

$ws.Range("H133").Value = 85000
$ws.Range("J133").Value = 85000
$ws.Range("L133").Value = 85000
$ws.Range("N133").Value = -95120

$ws.Range("H136").Value = 18370116
$ws.Range("I136").Value = 20524078
$ws.Range("J136").Value = 61440.5
$ws.Range("K136").Value = 61572234
$ws.Range("L136").Value = 184321.5
$ws.Range("M136").Value = -61569684
$ws.Range("N136").Value = -189421.5

$ws.Range("H141").Value = 250077420
$ws.Range("J141").Value = 250077420
$ws.Range("L141").Value = 250077420
$ws.Range("N141").Value = -250087780

